$p = $ppt.ActivePresentation

# Slide 1: add an "Appear" (With Previous) entrance effect to the
# source-credit textbox (shape id=4, "Textfeld 3").
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$null = $s1.TimeLine.MainSequence.AddEffect($shp1, 1, 0, 2)

# Slide 10: add the same "Appear" (With Previous) entrance effect to the
# source-credit textbox (shape id=10, "Textfeld 9").
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(4)
$null = $s10.TimeLine.MainSequence.AddEffect($shp10, 1, 0, 2)
